# Applies the scheduled-runner update to Anima_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2942.2666
$ws.Range("J76").Value = 2835.3333
$ws.Range("L76").Value = 2835.3333
$ws.Range("N76").Value = -3465.3333
$ws.Range("H79").Value = 2942.2666
$ws.Range("J79").Value = 2835.3333
$ws.Range("L79").Value = 2835.3333
$ws.Range("N79").Value = -5019.3333
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H125").Value = 1707.8889
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 1796.375
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 16167.375
$ws.Range("M125").Value = -6540
$ws.Range("N125").Value = -21087.375
$ws.Range("H138").Value = 1823.9231
$ws.Range("I138").Value = 1279.6904
$ws.Range("J138").Value = 2290.4082
$ws.Range("K138").Value = 3839.0712
$ws.Range("L138").Value = 6871.2246
$ws.Range("M138").Value = 1300.9288
$ws.Range("N138").Value = -17151.2246

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1005228.5
$ws.Range("I32").Value = 1075366.8
$ws.Range("J32").Value = 40828.5
$ws.Range("K32").Value = 1075366.8
$ws.Range("L32").Value = 40828.5
$ws.Range("M32").Value = -1075079.8
$ws.Range("N32").Value = -41402.5
$ws.Range("H45").Value = 3728
$ws.Range("I45").Value = 4299
$ws.Range("J45").Value = 2966.6667
$ws.Range("K45").Value = 4299
$ws.Range("L45").Value = 2966.6667
$ws.Range("M45").Value = -3922
$ws.Range("N45").Value = -3720.6667
$ws.Range("H61").Value = 3336.125
$ws.Range("I61").Value = 3288.5
$ws.Range("J61").Value = 3370.1428
$ws.Range("K61").Value = 3288.5
$ws.Range("L61").Value = 3370.1428
$ws.Range("M61").Value = -3076.5
$ws.Range("N61").Value = -3794.1428
$ws.Range("H88").Value = 2707
$ws.Range("J88").Value = 2707
$ws.Range("L88").Value = 2707
$ws.Range("N88").Value = -3519
$ws.Range("H91").Value = 2707
$ws.Range("J91").Value = 2707
$ws.Range("L91").Value = 2707
$ws.Range("N91").Value = -5515
$ws.Range("H136").Value = 3336.125
$ws.Range("I136").Value = 3288.5
$ws.Range("J136").Value = 3370.1428
$ws.Range("K136").Value = 9865.5
$ws.Range("L136").Value = 10110.4284
$ws.Range("M136").Value = -7315.5
$ws.Range("N136").Value = -15210.4284
$ws.Range("H137").Value = 27086.666
$ws.Range("J137").Value = 23347.5
$ws.Range("L137").Value = 23347.5
$ws.Range("N137").Value = -33547.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3180
$ws.Range("I86").Value = 3366.2
$ws.Range("J86").Value = 2249
$ws.Range("K86").Value = 3366.2
$ws.Range("L86").Value = 2249
$ws.Range("M86").Value = -2243.2
$ws.Range("N86").Value = -4495
$ws.Range("H89").Value = 3180
$ws.Range("I89").Value = 3366.2
$ws.Range("J89").Value = 2249
$ws.Range("K89").Value = 16831
$ws.Range("L89").Value = 11245
$ws.Range("M89").Value = -11215
$ws.Range("N89").Value = -22477
$ws.Range("H134").Value = 5471
$ws.Range("I134").Value = 5802.75
$ws.Range("J134").Value = 5205.6
$ws.Range("K134").Value = 17408.25
$ws.Range("L134").Value = 15616.8
$ws.Range("M134").Value = -14873.25
$ws.Range("N134").Value = -20686.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3114.647
$ws.Range("I31").Value = 983.551
$ws.Range("J31").Value = 6015.3057
$ws.Range("K31").Value = 983.551
$ws.Range("L31").Value = 6015.3057
$ws.Range("M31").Value = -688.551
$ws.Range("N31").Value = -6605.3057
$ws.Range("H34").Value = 3114.647
$ws.Range("I34").Value = 983.551
$ws.Range("J34").Value = 6015.3057
$ws.Range("K34").Value = 983.551
$ws.Range("L34").Value = 6015.3057
$ws.Range("M34").Value = -781.551
$ws.Range("N34").Value = -6419.3057
$ws.Range("H132").Value = 1848.3422
$ws.Range("I132").Value = 1669.9656
$ws.Range("J132").Value = 2423.111
$ws.Range("K132").Value = 5009.8968
$ws.Range("L132").Value = 7269.333
$ws.Range("M132").Value = -2479.8968
$ws.Range("N132").Value = -12329.333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3999.5833
$ws.Range("H79").Value = 3999.5833
$ws.Range("H120").Value = 10606
$ws.Range("J120").Value = 12500
$ws.Range("L120").Value = 37500
$ws.Range("N120").Value = -47176
$ws.Range("H132").Value = 2715.8286
$ws.Range("I132").Value = 2575.4
$ws.Range("J132").Value = 2772
$ws.Range("K132").Value = 23178.6
$ws.Range("L132").Value = 24948
$ws.Range("M132").Value = -20648.6
$ws.Range("N132").Value = -30008
$ws.Range("H137").Value = 7254011.5
$ws.Range("I137").Value = 27795760
$ws.Range("J137").Value = 3982.353
$ws.Range("K137").Value = 83387280
$ws.Range("L137").Value = 11947.059
$ws.Range("M137").Value = -83382180
$ws.Range("N137").Value = -22147.059
$ws.Range("H140").Value = 1631.875
$ws.Range("I140").Value = 1196
$ws.Range("K140").Value = 3588
$ws.Range("M140").Value = 1592
$ws.Range("H141").Value = 6678.1924
$ws.Range("I141").Value = 3402.5386
$ws.Range("J141").Value = 9953.846
$ws.Range("K141").Value = 10207.6158
$ws.Range("L141").Value = 29861.538
$ws.Range("M141").Value = -5027.6158
$ws.Range("N141").Value = -40221.538

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 24955.562
$ws.Range("I82").Value = 24298
$ws.Range("J82").Value = 24999.4
$ws.Range("K82").Value = 24298
$ws.Range("L82").Value = 24999.4
$ws.Range("M82").Value = -23915
$ws.Range("N82").Value = -25765.4
$ws.Range("H85").Value = 24955.562
$ws.Range("I85").Value = 24298
$ws.Range("J85").Value = 24999.4
$ws.Range("K85").Value = 24298
$ws.Range("L85").Value = 24999.4
$ws.Range("M85").Value = -22972
$ws.Range("N85").Value = -27651.4
$ws.Range("H113").Value = 102079.9
$ws.Range("I113").Value = 126849.875
$ws.Range("K113").Value = 126849.875
$ws.Range("M113").Value = -124679.875
$ws.Range("H138").Value = 49414.5
$ws.Range("J138").Value = 49414.5
$ws.Range("L138").Value = 49414.5
$ws.Range("N138").Value = -59694.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6433
$ws.Range("I7").Value = 7201.6
$ws.Range("K7").Value = 7201.6
$ws.Range("M7").Value = -7089.6
$ws.Range("H93").Value = 18826.334
$ws.Range("I93").Value = 50701.5
$ws.Range("J93").Value = 2888.75
$ws.Range("K93").Value = 50701.5
$ws.Range("L93").Value = 2888.75
$ws.Range("M93").Value = -49453.5
$ws.Range("N93").Value = -5384.75
$ws.Range("H126").Value = 6433
$ws.Range("I126").Value = 7201.6
$ws.Range("K126").Value = 21604.8
$ws.Range("M126").Value = -19134.8
$ws.Range("H132").Value = 2674.6428
$ws.Range("I132").Value = 2800.6
$ws.Range("J132").Value = 2529.3076
$ws.Range("K132").Value = 8401.799999999999
$ws.Range("L132").Value = 7587.9228
$ws.Range("M132").Value = -5871.799999999999
$ws.Range("N132").Value = -12647.9228
$ws.Range("H137").Value = 35750

Write-Host "Applied Anima_Profits updates"